# Adds 12 new rows (Power BI Service/Workspace/Dataflow/Gateway/Row-Level
# Security/Q&A, plus a new "Machine Learning" category block) to the
# "Tabela1" table on Arkusz1, growing it from A1:D161 to A1:D173, matching
# the author's upload diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Category, Concept, Definition, Example usage, row-height for each new row.
$rows = @(
    @("Power BI", "Power BI Service", "Usługa online do udostępniania raportów i pulpitów nawigacyjnych.", "Opublikowałem raport z Power BI Desktop do Power BI Service, aby inni mogli go przeglądać.", 28.8),
    @("Power BI", "Workspace", "Obszar roboczy w Power BI Service służący do organizacji raportów i zbiorów danych.", "Utworzono workspace 'Sprzedaż' do przechowywania raportów zespołu sprzedaży.", 28.8),
    @("Power BI", "Dataflow", "Procesy ETL definiowane w Power BI Service do przekształcania i ładowania danych.", "Dataflow pobiera dane z SQL i przygotowuje czystą tabelę do raportów.", 28.8),
    @("Power BI", "Gateway", "Usługa umożliwiająca odświeżanie danych lokalnych w Power BI Service.", "Zainstalowano gateway, aby automatycznie odświeżać dane z lokalnej bazy danych.", 28.8),
    @("Power BI", "Row-Level Security", "Mechanizm ograniczający dostęp do danych na poziomie wiersza w modelu.", "Skonfigurowano RLS, aby sprzedawcy widzieli tylko dane swojego regionu.", 28.8),
    @("Power BI", "Q&A", "Funkcja zadawania pytań w języku naturalnym, zwracająca wizualizacje.", "Użytkownik wpisał 'sprzedaż w marcu' i otrzymał wykres linii.", 28.8),
    @("Machine Learning", "Train/Test Split", "Podział zbioru danych na części treningową i testową w celu oceny modelu.", "80% danych użyto do trenowania, 20% do testowania modelu.", 28.8),
    @("Machine Learning", "Feature Scaling", "Normalizacja lub standaryzacja cech, aby ulepszyć efektywność algorytmu.", "Zastosowano standaryzację, aby zmienne miały średnią 0 i odchylenie standardowe 1.", 28.8),
    @("Machine Learning", "Regularization", "Technika zapobiegająca nadmiernemu dopasowaniu poprzez dodanie kary za złożoność modelu.", "Użyto L2 regularizacji w regresji liniowej.", 43.2),
    @("Machine Learning", "Decision Tree", "Algorytm tworzący model w postaci drzewa decyzyjnego, dzielącego dane według cech.", "Drzewo decyzyjne sklasyfikowało klientów jako potencjalnych nabywców.", 28.8),
    @("Machine Learning", "Random Forest", "Ensemble złożony z wielu drzew decyzyjnych w celu zwiększenia stabilności i dokładności.", "Random Forest osiągnął 92% dokładności na zbiorze testowym.", 28.8),
    @("Machine Learning", "KNN (K-Nearest Neighbors)", "Algorytm klasyfikacji oparty na najbliższych sąsiadach w przestrzeni cech.", "KNN z k=5 sklasyfikował nowy punkt na podstawie pięciu najbliższych przykładów.", 28.8)
)

$startRow = 162

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $prevRow = $r - 1

    # Carry the established look (bold category column, wrapped/centered
    # text, "Arial Unicode MS" example column, …) down from the row above
    # instead of re-deriving it property by property.
    $ws.Range("A" + $prevRow + ":D" + $prevRow).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]

    $ws.Rows.Item($r).RowHeight = $rows[$i][4]
}

$excel.CutCopyMode = $false

$endRow = $startRow + $rows.Count - 1

# Grow the table ("Tabela1") to cover the newly-added rows.
$tbl = $ws.ListObjects.Item("Tabela1")
$newTableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($endRow, 4))
$tbl.Resize($newTableRange)

# Leave the view roughly where the author left it.
$ws.Range("D171").Select()
